$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NIK (A) and No_HP (H) columns contain long numeric strings that Excel
# would otherwise auto-convert to Number (losing precision / leading zeros).
# Jumlah (E) contains plain numeric text ("0", "65000", ...) that must stay
# text too (row 22 is the one genuine numeric exception, handled below).
# Mark those columns as Text ("@") before writing so the literal digits survive.

# --- Fix E9: was stored as a real number 65000, must become text "65000" ---
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "65000"

# --- Append the new payment history rows 10-21 (all-text data) ---
$ws.Range("A10:A21").NumberFormat = "@"
$ws.Range("E10:E21").NumberFormat = "@"
$ws.Range("H10:H21").NumberFormat = "@"

# Row 10
$ws.Cells.Item(10, 1).Value = '1234456278949542'
$ws.Cells.Item(10, 2).Value = 'BG4576HI'
$ws.Cells.Item(10, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(10, 4).Value = '05-08-2025 03:16'
$ws.Cells.Item(10, 5).Value = '65000'
$ws.Cells.Item(10, 6).Value = 'GoPay'
$ws.Cells.Item(10, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(10, 8).Value = '085267947261'
$ws.Cells.Item(10, 9).Value = 'Palembang'
$ws.Cells.Item(10, 10).Value = 'JNE'

# Row 11
$ws.Cells.Item(11, 1).Value = '1234456278949542'
$ws.Cells.Item(11, 2).Value = 'BG4576HI'
$ws.Cells.Item(11, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(11, 4).Value = '05-08-2025 03:46'
$ws.Cells.Item(11, 5).Value = '50000'
$ws.Cells.Item(11, 6).Value = 'Bank Mandiri'
$ws.Cells.Item(11, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(11, 8).Value = '085267947261'
$ws.Cells.Item(11, 9).Value = 'Palembang'
$ws.Cells.Item(11, 10).Value = 'JNE'

# Row 12
$ws.Cells.Item(12, 1).Value = '1234456278949542'
$ws.Cells.Item(12, 2).Value = 'BG4576HI'
$ws.Cells.Item(12, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(12, 4).Value = '05-08-2025 03:52'
$ws.Cells.Item(12, 5).Value = '50000'
$ws.Cells.Item(12, 6).Value = 'Bank Mandiri'
$ws.Cells.Item(12, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(12, 8).Value = '085267947261'
$ws.Cells.Item(12, 9).Value = 'Palembang'
$ws.Cells.Item(12, 10).Value = 'JNE'

# Row 13
$ws.Cells.Item(13, 1).Value = '1234456278949542'
$ws.Cells.Item(13, 2).Value = 'BG4576HI'
$ws.Cells.Item(13, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(13, 4).Value = '05-08-2025 04:04'
$ws.Cells.Item(13, 5).Value = '50000'
$ws.Cells.Item(13, 6).Value = 'Bank Mandiri'
$ws.Cells.Item(13, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(13, 8).Value = '085267947261'
$ws.Cells.Item(13, 9).Value = 'Palembang'
$ws.Cells.Item(13, 10).Value = 'JNE'

# Row 14
$ws.Cells.Item(14, 1).Value = '1234456278949542'
$ws.Cells.Item(14, 2).Value = 'BG4576HI'
$ws.Cells.Item(14, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(14, 4).Value = '05-08-2025 04:05'
$ws.Cells.Item(14, 5).Value = '50000'
$ws.Cells.Item(14, 6).Value = 'Bank Mandiri'
$ws.Cells.Item(14, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(14, 8).Value = '085267947261'
$ws.Cells.Item(14, 9).Value = 'Palembang'
$ws.Cells.Item(14, 10).Value = 'JNE'

# Row 15
$ws.Cells.Item(15, 1).Value = '1234456278949542'
$ws.Cells.Item(15, 2).Value = 'BG4576HI'
$ws.Cells.Item(15, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(15, 4).Value = '05-08-2025 04:07'
$ws.Cells.Item(15, 5).Value = '50000'
$ws.Cells.Item(15, 6).Value = 'SeaBank'
$ws.Cells.Item(15, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(15, 8).Value = '085267947261'
$ws.Cells.Item(15, 9).Value = 'Palembang'
$ws.Cells.Item(15, 10).Value = 'JNE'

# Row 16
$ws.Cells.Item(16, 1).Value = '1234456278949542'
$ws.Cells.Item(16, 2).Value = 'BG4576HI'
$ws.Cells.Item(16, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(16, 4).Value = '05-08-2025 04:09'
$ws.Cells.Item(16, 5).Value = '50000'
$ws.Cells.Item(16, 6).Value = 'SeaBank'
$ws.Cells.Item(16, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(16, 8).Value = '085267947261'
$ws.Cells.Item(16, 9).Value = 'Palembang'
$ws.Cells.Item(16, 10).Value = 'JNE'

# Row 17
$ws.Cells.Item(17, 1).Value = '1234456278949542'
$ws.Cells.Item(17, 2).Value = 'BG4576HI'
$ws.Cells.Item(17, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(17, 4).Value = '05-08-2025 04:10'
$ws.Cells.Item(17, 5).Value = '50000'
$ws.Cells.Item(17, 6).Value = 'SeaBank'
$ws.Cells.Item(17, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(17, 8).Value = '085267947261'
$ws.Cells.Item(17, 9).Value = 'Palembang'
$ws.Cells.Item(17, 10).Value = 'JNE'

# Row 18
$ws.Cells.Item(18, 1).Value = '1234456278949542'
$ws.Cells.Item(18, 2).Value = 'BG4576HI'
$ws.Cells.Item(18, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(18, 4).Value = '05-08-2025 11:16'
$ws.Cells.Item(18, 5).Value = '50000'
$ws.Cells.Item(18, 6).Value = 'SeaBank'
$ws.Cells.Item(18, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(18, 8).Value = '085267947261'
$ws.Cells.Item(18, 9).Value = 'Palembang'
$ws.Cells.Item(18, 10).Value = 'JNE'

# Row 19
$ws.Cells.Item(19, 1).Value = '1234456278949542'
$ws.Cells.Item(19, 2).Value = 'BG4576HI'
$ws.Cells.Item(19, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(19, 4).Value = '05-08-2025 11:17'
$ws.Cells.Item(19, 5).Value = '0'
$ws.Cells.Item(19, 6).Value = 'Bank Rakyat Indonesia (BRI)'
$ws.Cells.Item(19, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(19, 8).Value = '085267947261'
$ws.Cells.Item(19, 9).Value = 'Palembang'
$ws.Cells.Item(19, 10).Value = 'JNE'

# Row 20
$ws.Cells.Item(20, 1).Value = '1234456278949542'
$ws.Cells.Item(20, 2).Value = 'BG4576HI'
$ws.Cells.Item(20, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(20, 4).Value = '08-05-2025 11:41'
$ws.Cells.Item(20, 5).Value = '0'
$ws.Cells.Item(20, 6).Value = 'Bank Central Asia (BCA)'
$ws.Cells.Item(20, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(20, 8).Value = '085267947261'
$ws.Cells.Item(20, 9).Value = 'Palembang'
$ws.Cells.Item(20, 10).Value = 'JNE'

# Row 21
$ws.Cells.Item(21, 1).Value = '1234456278949542'
$ws.Cells.Item(21, 2).Value = 'BG4576HI'
$ws.Cells.Item(21, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(21, 4).Value = '05-08-2025 11:46'
$ws.Cells.Item(21, 5).Value = '40000'
$ws.Cells.Item(21, 6).Value = 'Bank Rakyat Indonesia (BRI)'
$ws.Cells.Item(21, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(21, 8).Value = '085267947261'
$ws.Cells.Item(21, 9).Value = 'Palembang'
$ws.Cells.Item(21, 10).Value = 'JNE'

# --- Row 22: same as above, but Jumlah (E22) is a genuine number (30000) ---
$ws.Range("A22").NumberFormat = "@"
$ws.Range("H22").NumberFormat = "@"
$ws.Cells.Item(22, 1).Value = '1234456278949542'
$ws.Cells.Item(22, 2).Value = 'BG4576HI'
$ws.Cells.Item(22, 3).Value = 'Nia Rahmadani'
$ws.Cells.Item(22, 4).Value = '05-08-2025 11:50'
$ws.Cells.Item(22, 5).Value = 30000
$ws.Cells.Item(22, 6).Value = 'Bank Rakyat Indonesia (BRI)'
$ws.Cells.Item(22, 7).Value = 'Nia Rahmadani'
$ws.Cells.Item(22, 8).Value = '085267947261'
$ws.Cells.Item(22, 9).Value = 'Palembang'
$ws.Cells.Item(22, 10).Value = 'JNE'
